# issue #5: property land done
# Completes the 土地 (Land) sheet with the canonical field-name header row
# and the extra metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) that the other
# sheets already carry. Also cleans up stray whitespace / stray dashes
# that had crept into a handful of shared strings across sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 土地 (Land) -- header rename + 7 new trailing columns (I..O)
# ---------------------------------------------------------------------
$land = $wb.Worksheets.Item("土地")

function Set-HeaderCell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $text
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160    # xlTop
    $c.Borders.LineStyle = 1
}

# existing header cells B1:H1 get renamed from Chinese labels to the
# canonical snake_case field names
Set-HeaderCell $land 1 2 "name"
Set-HeaderCell $land 1 3 "area"
Set-HeaderCell $land 1 4 "share_portion"
Set-HeaderCell $land 1 5 "owner"
Set-HeaderCell $land 1 6 "register_date"
Set-HeaderCell $land 1 7 "register_reason"
Set-HeaderCell $land 1 8 "acquire_value"

# new trailing header cells I1:O1
Set-HeaderCell $land 1 9  "property_category"
Set-HeaderCell $land 1 10 "category"
Set-HeaderCell $land 1 11 "date"
Set-HeaderCell $land 1 12 "legislator_name"
Set-HeaderCell $land 1 13 "legislator_id"
Set-HeaderCell $land 1 14 "source_file"
Set-HeaderCell $land 1 15 "index"

# new data columns I..O for the 3 existing data rows (2,3,4)
$land.Cells.Item(2, 9).Value  = "land"
$land.Cells.Item(2, 10).Value = "normal"
$land.Cells.Item(2, 11).Value = "2012-03-26"
$land.Cells.Item(2, 12).Value = "潘維剛"
$land.Cells.Item(2, 13).Value = 678
$land.Cells.Item(2, 14).Value = "tmp71a01"
$land.Cells.Item(2, 15).Value = 14

$land.Cells.Item(3, 9).Value  = "land"
$land.Cells.Item(3, 10).Value = "normal"
$land.Cells.Item(3, 11).Value = "2012-03-26"
$land.Cells.Item(3, 12).Value = "潘維剛"
$land.Cells.Item(3, 13).Value = 678
$land.Cells.Item(3, 14).Value = "tmp71a01"
$land.Cells.Item(3, 15).Value = 15

$land.Cells.Item(4, 9).Value  = "land"
$land.Cells.Item(4, 10).Value = "normal"
$land.Cells.Item(4, 11).Value = "2012-03-26"
$land.Cells.Item(4, 12).Value = "潘維剛"
$land.Cells.Item(4, 13).Value = 678
$land.Cells.Item(4, 14).Value = "tmp71a01"
$land.Cells.Item(4, 15).Value = 16

# clean up stray whitespace inside a couple of the land-parcel / date strings
$land.Cells.Item(2, 2).Value = "新北市八里區小八里分段楓櫃斗湖小段02690003地號"
$land.Cells.Item(2, 6).Value = "101年01月11曰"
$land.Cells.Item(3, 2).Value = "臺北市松山區敦化段五小段00290010地號"
$land.Cells.Item(3, 6).Value = "72年06月11曰"
$land.Cells.Item(4, 2).Value = "新北市八里區小八里分段楓櫃斗湖小段02830001地號"
$land.Cells.Item(4, 6).Value = "101年03月08日"

$land.Range("A1:O4").EntireColumn.AutoFit()

# ---------------------------------------------------------------------
# 建物 (Building) -- whitespace / dash cleanup only
# ---------------------------------------------------------------------
$building = $wb.Worksheets.Item("建物")
$building.Cells.Item(2, 2).Value = "臺北市松山區敦化段五小段02392000建號"
$building.Cells.Item(2, 6).Value = "72年06月11曰"

# ---------------------------------------------------------------------
# 汽車 (Car) -- whitespace cleanup only
# ---------------------------------------------------------------------
$car = $wb.Worksheets.Item("汽車")
$car.Cells.Item(2, 5).Value = "97年07月07日"

# ---------------------------------------------------------------------
# 存款 (Deposit) -- whitespace cleanup only
# ---------------------------------------------------------------------
$deposit = $wb.Worksheets.Item("存款")
$deposit.Cells.Item(2, 2).Value  = "立法院郵局（第25支局）"
$deposit.Cells.Item(3, 2).Value  = "台北長安郵局（第46支局）"
$deposit.Cells.Item(6, 2).Value  = "日盛國際商業銀行松山分行"
$deposit.Cells.Item(8, 2).Value  = "中國信託商業銀行城中分行"
$deposit.Cells.Item(10, 2).Value = "台新國際商業銀行敦北分行"
$deposit.Cells.Item(16, 2).Value = "台北富邦商業銀行敦南分行"
$deposit.Cells.Item(20, 2).Value = "國泰世華商業銀行南京東路分行"

# ---------------------------------------------------------------------
# 基金受益憑證 (Fund) -- whitespace cleanup only
# ---------------------------------------------------------------------
$fund = $wb.Worksheets.Item("基金受益憑證")
$fund.Cells.Item(3, 2).Value = "聯博全球高收益債券AT股"

# ---------------------------------------------------------------------
# 保險 (Insurance) -- whitespace cleanup only
# ---------------------------------------------------------------------
$insurance = $wb.Worksheets.Item("保險")
$insurance.Cells.Item(2, 3).Value = "吉祥變額萬能終身壽險(A型）"
$insurance.Cells.Item(3, 3).Value = "吉祥變額萬能終身壽險(A型）"
$insurance.Cells.Item(4, 3).Value = "限期繳費單利增值終身壽險已型"

# ---------------------------------------------------------------------
# 債務 (Debt) -- whitespace cleanup only
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item("債務")
$debt.Cells.Item(2, 4).Value = "永豐銀行西松分行臺北市松山區東興路"
$debt.Cells.Item(2, 6).Value = "95年10月27日"
$debt.Cells.Item(3, 4).Value = "永豐銀行板新分行新北市板橋區民權路"
$debt.Cells.Item(3, 6).Value = "96年07月25F1"

# ---------------------------------------------------------------------
# 事業投資 (Business investment) -- whitespace cleanup only
# ---------------------------------------------------------------------
$invest = $wb.Worksheets.Item("事業投資")
$invest.Cells.Item(2, 4).Value = "臺北市杭州南路1段63號5樓之1"
$invest.Cells.Item(2, 6).Value = "98年08月10日"
$invest.Cells.Item(3, 3).Value = "傳智國際文化事業股份有限公司"
$invest.Cells.Item(3, 4).Value = "臺北市羅斯福路2段116號3樓"
$invest.Cells.Item(3, 6).Value = "89年08月11曰"
$invest.Cells.Item(4, 4).Value = "臺北市民權東路3段106巷36號7樓"
$invest.Cells.Item(4, 6).Value = "92年03月20日"
$invest.Cells.Item(5, 4).Value = "新北市寶僑路235巷6弄5號5樓"
$invest.Cells.Item(5, 6).Value = "92年03月20H"
